# Apply updated crypto price (D) and volume% (E) values for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''22.063.60'
$ws.Range("E2").Value = '  -1.51%  '
$ws.Range("D3").Value = '''1.557.54'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''0.9998'
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").Value = '''287.64'
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").Value = '''0.3879'
$ws.Range("E7").Value = '  +3.99%  '
$ws.Range("D8").Value = '''0.3229'
$ws.Range("E8").Value = '  -1.65%  '
$ws.Range("D9").Value = '''42.78'
$ws.Range("E9").Value = '  -6.13%  '
$ws.Range("D10").Value = '''1.121'
$ws.Range("E10").Value = '  -2.03%  '
$ws.Range("D11").Value = '''0.07357'
$ws.Range("E11").Value = '  -0.85%  '
$ws.Range("D12").Value = '''1.001'
$ws.Range("E12").Value = '  -0.01%  '
$ws.Range("D13").Value = '''19.30'
$ws.Range("E13").Value = '  -5.74%  '
$ws.Range("D14").Value = '''5.696'
$ws.Range("E14").Value = '  -2.66%  '
$ws.Range("D15").Value = '''6.798'
$ws.Range("E15").Value = '  -0.69%  '
$ws.Range("E16").Value = '  +2.44%  '
$ws.Range("D17").Value = '''1.553.14'
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("D18").Value = '''0.06615'
$ws.Range("E18").Value = '  -1.10%  '
$ws.Range("D19").Value = '''85.24'
$ws.Range("E19").Value = '  -0.77%  '
$ws.Range("D20").Value = '''6.397'
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("D21").Value = '''0.9994'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").Value = '''15.97'
$ws.Range("E22").Value = '  -2.00%  '
$ws.Range("E23").Value = '  -2.44%  '
$ws.Range("D24").Value = '''22.077.32'
$ws.Range("E24").Value = '  -1.41%  '
$ws.Range("D25").Value = '''2.337'
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("D26").Value = '''2.560'
$ws.Range("E26").Value = '  -0.64%  '
$ws.Range("D27").Value = '''148.05'
$ws.Range("E27").Value = '  -2.61%  '
$ws.Range("E28").Value = '  -2.31%  '
$ws.Range("D29").Value = '''4.860'
$ws.Range("E29").Value = '  -1.25%  '
$ws.Range("D30").Value = '''1.728.39'
$ws.Range("E30").Value = '  -0.92%  '
$ws.Range("D31").Value = '''120.94'
$ws.Range("E31").Value = '  -2.22%  '
$ws.Range("D32").Value = '''1.100'
$ws.Range("E32").Value = '  +2.32%  '
$ws.Range("D33").Value = '''5.835'
$ws.Range("E33").Value = '  -2.08%  '
$ws.Range("D34").Value = '''1.684'
$ws.Range("E34").Value = '  -13.79%  '
$ws.Range("D35").Value = '''9.416'
$ws.Range("E35").Value = '  -3.07%  '
$ws.Range("E36").Value = '  -0.70%  '
$ws.Range("D37").Value = '''0.06261'
$ws.Range("E37").Value = '  -0.59%  '
$ws.Range("D38").Value = '''0.02305'
$ws.Range("E38").Value = '  -3.31%  '
$ws.Range("D39").Value = '''5.245'
$ws.Range("E39").Value = '  -0.92%  '
$ws.Range("D40").Value = '''0.2119'
$ws.Range("E40").Value = '  -3.27%  '
$ws.Range("D41").Value = '''1.225'
$ws.Range("E41").Value = '  -5.90%  '
$ws.Range("D42").Value = '''10.88'
$ws.Range("E42").Value = '  -2.32%  '
$ws.Range("D43").Value = '''0.9991'
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").Value = '''0.5964'
$ws.Range("E44").Value = '  -2.44%  '
$ws.Range("D45").Value = '''13.58'
$ws.Range("E45").Value = '  -1.61%  '
$ws.Range("D46").Value = '''3.719'
$ws.Range("E46").Value = '  -0.81%  '
$ws.Range("D47").Value = '''0.5751'
$ws.Range("E47").Value = '  -3.35%  '
$ws.Range("D48").Value = '''1.934'
$ws.Range("E48").Value = '  -3.70%  '
$ws.Range("D49").Value = '''119.36'
$ws.Range("E49").Value = '  -3.88%  '
$ws.Range("E50").Value = '  -1.60%  '
$ws.Range("D51").Value = '''0.06892'
$ws.Range("E51").Value = '  -3.66%  '
